# Applies the statistics_summary edits:
#  - Rename header B1 "BTC" -> "IBTC"
#  - Fill in previously-empty Mean/Median/StdDev/Kurtosis/Skewness cells for
#    columns B, D, E and update existing numeric values across B:F for rows 2-9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("B1").Value = "IBTC"

# Row 2 - Mean
$ws.Range("B2").Value = 34.02081627210884
$ws.Range("C2").Value = 127376.8027210884
$ws.Range("D2").Value = 9875.46084602721
$ws.Range("E2").Value = 7.855714285714286
$ws.Range("F2").Value = 285.7742180816327

# Row 3 - Median
$ws.Range("B3").Value = 36.209999
$ws.Range("C3").Value = 127668
$ws.Range("D3").Value = 9886.252229
$ws.Range("E3").Value = 7.91
$ws.Range("F3").Value = 288.290009

# Row 4 - Standard Deviation
$ws.Range("B4").Value = 5.842908680803125
$ws.Range("C4").Value = 2017.562372534595
$ws.Range("D4").Value = 38.76511435206717
$ws.Range("E4").Value = 0.2591115117284527
$ws.Range("F4").Value = 12.97404173256524

# Row 5 - Kurtosis
$ws.Range("B5").Value = -1.017277626065159
$ws.Range("C5").Value = 0.8240942329896561
$ws.Range("D5").Value = -0.2770808927941548
$ws.Range("E5").Value = -0.4171455961159491
$ws.Range("F5").Value = -0.4368381318723622

# Row 6 - Skewness
$ws.Range("B6").Value = -0.6936905131393687
$ws.Range("C6").Value = -0.805890693663089
$ws.Range("D6").Value = -0.6165237100922017
$ws.Range("E6").Value = -0.5060570587828281
$ws.Range("F6").Value = -0.2494028448025493

# Row 7 - Fishers Information
$ws.Range("B7").Value = 0.06199446552504129
$ws.Range("C7").Value = 0.09843999638032346
$ws.Range("D7").Value = 0.09308118403278205
$ws.Range("E7").Value = 0.0980346371114356
$ws.Range("F7").Value = 0.1000914157361388

# Row 8 - MIEE
$ws.Range("B8").Value = 0.1869231013248247
$ws.Range("C8").Value = 0.236840330337551
$ws.Range("D8").Value = 0.2260385773870736
$ws.Range("E8").Value = 0.2342371191902062
$ws.Range("F8").Value = 0.2965468887295143

# Row 9 - Permutation Entropy
$ws.Range("B9").Value = 0.8750713642002167
$ws.Range("C9").Value = 0.8615996660427724
$ws.Range("D9").Value = 0.8670426066457084
$ws.Range("E9").Value = 0.8637975179212295
$ws.Range("F9").Value = 0.8035445270066246
